$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "ECs"

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8467519999999999
$ws.Range("H2").Value = 2.540256
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.438062
$ws.Range("N2").Value = 1.314186
$ws.Range("O2").Value = 0.6074000808827777
$ws.Range("P2").Value = 0.6074000808827777
$ws.Range("Q2").Value = 0.370929874624
$ws.Range("R2").Value = 3.338368871616
$ws.Range("S2").Value = 0.6074000808827777
$ws.Range("T2").Value = 0.6074000808827777

$ws.Range("A3").Value = "M2"
$ws.Range("B3").Value = "Ccl24"
$ws.Range("C3").Value = "Ccr3"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8467519999999999
$ws.Range("H3").Value = 2.540256
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2831463333333333
$ws.Range("N3").Value = 0.8494390000000001
$ws.Range("O3").Value = 0.3925999191172223
$ws.Range("P3").Value = 0.3925999191172223
$ws.Range("Q3").Value = 0.2397547240426667
$ws.Range("R3").Value = 2.157792516384
$ws.Range("S3").Value = 0.3925999191172223
$ws.Range("T3").Value = 0.3925999191172223
